# Add "schmitt trigger" (SN74ACT14) to the parts list.
#
# Layout change: a new "Info" column is inserted between "Part #" (B) and
# "Datasheet" (which moves from C to D); the old, never-populated "Qty"
# column (old D) is removed. A new row for the SN74ACT14 part is inserted
# between the "not / SN74HCS05" row and the "dmx 3 bit / SN74HCS137" row,
# and the dmx row also gets an "out of stock" Info note.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Structural edits -------------------------------------------------
# New "Info" column: inherits column B's formatting/width.
$ws.Columns.Item(3).Insert()
# Old "Qty" column (header only, never used) is now column E - drop it.
$ws.Columns.Item(5).Delete()
# Room for the new SN74ACT14 row, between current rows 5 and 6.
$ws.Rows.Item(6).Insert()

# The hyperlinks below survived the column insert/delete pointing at their
# old (now wrong) column, so drop them and recreate them against column D.
$ws.Hyperlinks.Delete()

# --- New part row (SN74ACT14 / schmitt trigger) ------------------------
$ws.Range("B6").Value2 = "SN74ACT14"
$ws.Range("D6").Value2 = "https://www.ti.com/lit/ds/symlink/sn74act14.pdf?HQS=dis-mous-null-mousermode-dsf-pf-null-wwe&ts=1640459360692&ref_url=https%253A%252F%252Fwww.mouser.ca%252F"

# --- Info column header + values ---------------------------------------
$ws.Range("C1").Value2 = "Info"
$ws.Range("C1").Font.Bold = $true

$ws.Range("C6").Value2 = "schmitt trigger"
$ws.Range("C7").Value2 = "out of stock"

# --- Recreate hyperlinks against column D -------------------------------
function Add-DatasheetLink($addr, $url, $fill) {
    $ws.Hyperlinks.Add($ws.Range($addr), $url) | Out-Null
    # Hyperlinks.Add() stamps its own style variant; reapply the sheet's
    # normal Hyperlink look (and the column's "fill" alignment where the
    # rest of the table uses it) on top of it.
    $ws.Range($addr).Style = "Hyperlink"
    if ($fill) {
        $ws.Range($addr).HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignFill
    }
}

Add-DatasheetLink "D2" "https://www.ti.com/lit/ds/symlink/sn74hcs86.pdf?ts=1638720224397&ref_url=https%253A%252F%252Fwww.ti.com%252Flogic-voltage-translation%252Flogic-gates%252Fexclusive-or-xor-gates%252Fproducts.html" $true
Add-DatasheetLink "D3" "https://www.ti.com/lit/ds/symlink/sn74act08.pdf?ts=1638720379375&ref_url=https%253A%252F%252Fwww.ti.com%252Fsitesearch%252Fdocs%252Funiversalsearch.tsp%253FlangPref%253Den-US%2526searchTerm%253DSN74ACT08%2526nr%253D22" $true
Add-DatasheetLink "D4" "https://www.ti.com/lit/ds/symlink/sn74hcs32.pdf?ts=1638697246303" $true
Add-DatasheetLink "D5" "https://www.ti.com/lit/gpn/sn74hcs05" $true
Add-DatasheetLink "D6" "https://www.ti.com/lit/ds/symlink/sn74act14.pdf?HQS=dis-mous-null-mousermode-dsf-pf-null-wwe&ts=1640459360692&ref_url=https%253A%252F%252Fwww.mouser.ca%252F" $true
Add-DatasheetLink "D7" "https://www.ti.com/lit/gpn/sn74hcs137-q1" $true
Add-DatasheetLink "D8" "https://datasheet.lcsc.com/lcsc/1809191939_Nexperia-74HC138PW-118_C47455.pdf" $false
Add-DatasheetLink "D9" "https://datasheet.lcsc.com/lcsc/1912141436_ISSI-Integrated-Silicon-Solution-IS61C5128AS-25QLI_C443416.pdf" $true

# --- View state ----------------------------------------------------------
$ws.Range("B14").Select()

Write-Output "done"
